# Cork KPI crawler - add "reduced emissions for cars" (5% / 10%) columns
# Inserts two new columns (D, E) between "Emissions for cars" (C) and
# "Emissions for boats" (old D, now F), each holding a formula that derives
# a reduced-emissions scenario from the "Emissions for cars" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (old D,E,F shift right to F,G,H)
$ws.Range("D:E").Insert()

# New columns should keep the same width as column C (they hold the same
# kind of numeric data)
$ws.Range("D:E").ColumnWidth = $ws.Range("C1").ColumnWidth

# Header row - new shared strings
$ws.Range("D1").Value = "Reduced Emissions for cars (5%)"
$ws.Range("E1").Value = "Reduced Emissions for cars (10%)"

# Header row now wraps onto two lines given the extra columns
$ws.Rows.Item(1).RowHeight = 58

# Row 2 formulas (entered individually)
$ws.Range("D2").Formula = "=0.95*C2"
$ws.Range("E2").Formula = "=0.9*C2"

# Rows 3-7 formulas (entered/filled as one block)
$ws.Range("D3:D7").Formula = "=0.95*C3"
$ws.Range("E3:E7").Formula = "=0.9*C3"

# Restore the last active selection
$ws.Range("E10").Select() | Out-Null
